$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.591
$ws.Range("C3").Value = -12.634
$ws.Range("C5").Value = -12.261
$ws.Range("E7").Value = 12.752
$ws.Range("A9").Value = -20.775
$ws.Range("E9").Value = 12.716
$ws.Range("C11").Value = -12.917
$ws.Range("C12").Value = -12.628
$ws.Range("A13").Value = -21.99
$ws.Range("A16").Value = -20.866
$ws.Range("A18").Value = -21.577
$ws.Range("A20").Value = -21.682
$ws.Range("C21").Value = -12.816
$ws.Range("E21").Value = 13.172
